# Study of the denomination of the credits
# Renames the "MH" column to "Denomination" and replaces its Yes/No
# values with the new Notable / Outstanding / Honors denominations
# (derived from the Mark column), then marks a new working cell (G46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header: column E used to be "MH" -> now "Denomination"
$ws.Range("E1").Value = "Denomination"

# 2) Column E data rows 2-32: replace with the new denomination values
#    (rows 33-39 keep their existing "-" placeholder and are untouched)
$denominations = @{
    2  = "Notable"
    3  = "Notable"
    4  = "Notable"
    5  = "Notable"
    6  = "Notable"
    7  = "Notable"
    8  = "Notable"
    9  = "Outstanding"
    10 = "Outstanding"
    11 = "Honors"
    12 = "Honors"
    13 = "Notable"
    14 = "Outstanding"
    15 = "Notable"
    16 = "Outstanding"
    17 = "Notable"
    18 = "Outstanding"
    19 = "Notable"
    20 = "Outstanding"
    21 = "Notable"
    22 = "Outstanding"
    23 = "Notable"
    24 = "Notable"
    25 = "Notable"
    26 = "Outstanding"
    27 = "Notable"
    28 = "Honors"
    29 = "Notable"
    30 = "Honors"
    31 = "Notable"
    32 = "Outstanding"
}

foreach ($row in $denominations.Keys) {
    $ws.Cells.Item($row, 5).Value = $denominations[$row]
}

# 3) New working cell left behind further down the sheet, underlined
#    like the stray L34 cell already on the sheet.
$ws.Range("G46").Font.Underline = 2

# 4) Leave the selection where the author left it
$ws.Range("G46").Select()
